$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 101, pushing existing rows 101-110 down to 102-111.
$ws.Rows("101:101").Insert()

# Populate the newly inserted row 101 with the new record's data.
$ws.Cells.Item(101, 1).Value = 8
$ws.Cells.Item(101, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(101, 3).Value = "Coquimbo"
$ws.Cells.Item(101, 4).Value = 44984
$ws.Cells.Item(101, 5).Value = 4
$ws.Cells.Item(101, 6).Value = 100114007
$ws.Cells.Item(101, 7).Value = "Jengibre"
$ws.Cells.Item(101, 8).Value = "Sin especificar"
$ws.Cells.Item(101, 9).Value = "Primera"
$ws.Cells.Item(101, 10).Value = 400
$ws.Cells.Item(101, 11).Value = 18000
$ws.Cells.Item(101, 12).Value = 19000
$ws.Cells.Item(101, 13).Value = 18500
$ws.Cells.Item(101, 14).Value = "$/caja 13 kilos"
$ws.Cells.Item(101, 15).Value = "Perú"
$ws.Cells.Item(101, 16).Value = 1423
$ws.Cells.Item(101, 17).Value = 13
$ws.Cells.Item(101, 18).Value = "Hortaliza"
